$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 53.75
$ws.Range("I8").Value = 53.75
$ws.Range("K8").Value = 161.25
$ws.Range("M8").Value = -22.25
$ws.Range("H17").Value = 10914.728
$ws.Range("J17").Value = 11848.9
$ws.Range("L17").Value = 35546.7
$ws.Range("N17").Value = -35882.7
$ws.Range("H33").Value = 1127.3334
$ws.Range("I33").Value = 1237.625
$ws.Range("J33").Value = 245
$ws.Range("K33").Value = 1237.625
$ws.Range("L33").Value = 245
$ws.Range("M33").Value = -1008.625
$ws.Range("N33").Value = -703
$ws.Range("H34").Value = 11390.429
$ws.Range("I34").Value = 11390.429
$ws.Range("K34").Value = 11390.429
$ws.Range("M34").Value = -11187.429
$ws.Range("H36").Value = 11390.429
$ws.Range("I36").Value = 11390.429
$ws.Range("K36").Value = 11390.429
$ws.Range("M36").Value = -10675.429
$ws.Range("H40").Value = 1799.4642
$ws.Range("I40").Value = 1628.8235
$ws.Range("K40").Value = 1628.8235
$ws.Range("M40").Value = -1453.8235
$ws.Range("H86").Value = 4566.3335
$ws.Range("J86").Value = 4566.3335
$ws.Range("L86").Value = 4566.3335
$ws.Range("N86").Value = -6812.3335
$ws.Range("H87").Value = 74000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 4566.3335
$ws.Range("J89").Value = 4566.3335
$ws.Range("L89").Value = 22831.6675
$ws.Range("N89").Value = -34063.6675
$ws.Range("H90").Value = 74000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H98").Value = 1765.65
$ws.Range("I98").Value = 412.2
$ws.Range("J98").Value = 3119.1
$ws.Range("K98").Value = 412.2
$ws.Range("L98").Value = 3119.1
$ws.Range("M98").Value = 1085.8
$ws.Range("N98").Value = -6115.1
$ws.Range("H101").Value = 14286938
$ws.Range("I101").Value = 25000492
$ws.Range("K101").Value = 75001476
$ws.Range("M101").Value = -74999854
$ws.Range("H122").Value = 1765.65
$ws.Range("I122").Value = 412.2
$ws.Range("J122").Value = 3119.1
$ws.Range("K122").Value = 1236.6
$ws.Range("L122").Value = 9357.299999999999
$ws.Range("M122").Value = 1213.4
$ws.Range("N122").Value = -14257.3
$ws.Range("H135").Value = 150.75
$ws.Range("J135").Value = 536
$ws.Range("L135").Value = 4824
$ws.Range("N135").Value = -9894
$ws.Range("H137").Value = 2017.6957
$ws.Range("I137").Value = 1886.5
$ws.Range("K137").Value = 5659.5
$ws.Range("M137").Value = -3109.5
$ws.Range("H141").Value = 1460.2273
$ws.Range("I141").Value = 1460.2273
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4380.6819
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 799.3181000000004
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 676.1111
$ws.Range("I25").Value = 135.625
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 135.625
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 266.375
$ws.Range("N25").Value = -5804
$ws.Range("H32").Value = 2542.8386
$ws.Range("I32").Value = 2542.8386
$ws.Range("K32").Value = 2542.8386
$ws.Range("M32").Value = -2255.8386
$ws.Range("H61").Value = 2743.125
$ws.Range("I61").Value = 2743.125
$ws.Range("K61").Value = 2743.125
$ws.Range("M61").Value = -2531.125
$ws.Range("H74").Value = 2031
$ws.Range("I74").Value = 1978.9445
$ws.Range("K74").Value = 1978.9445
$ws.Range("M74").Value = -1104.9445
$ws.Range("H77").Value = 2031
$ws.Range("I77").Value = 1978.9445
$ws.Range("K77").Value = 9894.7225
$ws.Range("M77").Value = -5526.7225
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H136").Value = 2743.125
$ws.Range("I136").Value = 2743.125
$ws.Range("K136").Value = 8229.375
$ws.Range("M136").Value = -5679.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1297.5454
$ws.Range("I86").Value = 1081.8572
$ws.Range("J86").Value = 1675
$ws.Range("K86").Value = 1081.8572
$ws.Range("L86").Value = 1675
$ws.Range("M86").Value = 41.14280000000008
$ws.Range("N86").Value = -3921
$ws.Range("H89").Value = 1297.5454
$ws.Range("I89").Value = 1081.8572
$ws.Range("J89").Value = 1675
$ws.Range("K89").Value = 5409.286
$ws.Range("L89").Value = 8375
$ws.Range("M89").Value = 206.7139999999999
$ws.Range("N89").Value = -19607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4867.1113
$ws.Range("I16").Value = 4515
$ws.Range("J16").Value = 6099.5
$ws.Range("K16").Value = 4515
$ws.Range("L16").Value = 6099.5
$ws.Range("M16").Value = -4228
$ws.Range("N16").Value = -6673.5
$ws.Range("H31").Value = 3196.8
$ws.Range("I31").Value = 2367.6667
$ws.Range("K31").Value = 2367.6667
$ws.Range("M31").Value = -2072.6667
$ws.Range("H34").Value = 3196.8
$ws.Range("I34").Value = 2367.6667
$ws.Range("K34").Value = 2367.6667
$ws.Range("M34").Value = -2165.6667
$ws.Range("H58").Value = 1146.5
$ws.Range("J58").Value = 770.7143
$ws.Range("L58").Value = 770.7143
$ws.Range("N58").Value = -1176.7143
$ws.Range("H113").Value = 4867.1113
$ws.Range("I113").Value = 4515
$ws.Range("J113").Value = 6099.5
$ws.Range("K113").Value = 4515
$ws.Range("L113").Value = 6099.5
$ws.Range("M113").Value = -2345
$ws.Range("N113").Value = -10439.5
$ws.Range("H134").Value = 869
$ws.Range("I134").Value = 869
$ws.Range("K134").Value = 2607
$ws.Range("M134").Value = -72
$ws.Range("H136").Value = 1146.5
$ws.Range("J136").Value = 770.7143
$ws.Range("L136").Value = 2312.1429
$ws.Range("N136").Value = -7412.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22001020
$ws.Range("I4").Value = 55000000
$ws.Range("K4").Value = 165000000
$ws.Range("M4").Value = -164999888
$ws.Range("H5").Value = 957.3077
$ws.Range("J5").Value = 1975.3334
$ws.Range("L5").Value = 5926.0002
$ws.Range("N5").Value = -6150.0002
$ws.Range("H40").Value = 112.84615
$ws.Range("I40").Value = 39.75
$ws.Range("K40").Value = 159
$ws.Range("M40").Value = -90
$ws.Range("H97").Value = 500
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1004
$ws.Range("N97").ClearContents()
$ws.Range("H103").Value = 713.38464
$ws.Range("I103").Value = 248
$ws.Range("K103").Value = 744
$ws.Range("M103").Value = 135
$ws.Range("H104").Value = 2259.2
$ws.Range("I104").Value = 2259.2
$ws.Range("K104").Value = 6777.599999999999
$ws.Range("M104").Value = -4156.599999999999
$ws.Range("H113").Value = 1532.1177
$ws.Range("I113").Value = 1258
$ws.Range("J113").Value = 1616.4615
$ws.Range("K113").Value = 3774
$ws.Range("L113").Value = 4849.3845
$ws.Range("M113").Value = -1604
$ws.Range("N113").Value = -9189.3845
$ws.Range("H122").Value = 1000.6667
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 1002
$ws.Range("K122").Value = 8982
$ws.Range("L122").Value = 9018
$ws.Range("M122").Value = -6532
$ws.Range("N122").Value = -13918
$ws.Range("H132").Value = 7395.3335
$ws.Range("I132").Value = 5043.5
$ws.Range("K132").Value = 45391.5
$ws.Range("M132").Value = -42861.5
$ws.Range("H135").Value = 957.3077
$ws.Range("J135").Value = 1975.3334
$ws.Range("L135").Value = 17778.0006
$ws.Range("N135").Value = -22848.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 35.4
$ws.Range("I2").Value = 50.833332
$ws.Range("J2").Value = 12.25
$ws.Range("K2").Value = 50.833332
$ws.Range("L2").Value = 12.25
$ws.Range("M2").Value = 62.166668
$ws.Range("N2").Value = -238.25
$ws.Range("H28").Value = 904
$ws.Range("J28").Value = 904
$ws.Range("L28").Value = 904
$ws.Range("N28").Value = -1288
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1210
$ws.Range("H31").Value = 299.5
$ws.Range("I31").Value = 299.5
$ws.Range("K31").Value = 299.5
$ws.Range("M31").Value = -7.5
$ws.Range("H37").Value = 299.5
$ws.Range("I37").Value = 299.5
$ws.Range("K37").Value = 299.5
$ws.Range("M37").Value = -22.5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H102").Value = 917.06665
$ws.Range("I102").Value = 917.06665
$ws.Range("K102").Value = 917.06665
$ws.Range("M102").Value = 704.93335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 171333.33
$ws.Range("I4").Value = 170333.33
$ws.Range("J4").Value = 173333.33
$ws.Range("K4").Value = 170333.33
$ws.Range("L4").Value = 173333.33
$ws.Range("M4").Value = -170220.33
$ws.Range("N4").Value = -173559.33
$ws.Range("H5").Value = 265004.5
$ws.Range("I5").Value = 30009
$ws.Range("K5").Value = 30009
$ws.Range("M5").Value = -29896
$ws.Range("H22").Value = 1992.9286
$ws.Range("I22").Value = 1966.75
$ws.Range("J22").Value = 2150
$ws.Range("K22").Value = 1966.75
$ws.Range("L22").Value = 2150
$ws.Range("M22").Value = -1671.75
$ws.Range("N22").Value = -2740
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H27").Value = 1992.9286
$ws.Range("I27").Value = 1966.75
$ws.Range("J27").Value = 2150
$ws.Range("K27").Value = 1966.75
$ws.Range("L27").Value = 2150
$ws.Range("M27").Value = -1859.75
$ws.Range("N27").Value = -2364
$ws.Range("H28").Value = 171333.33
$ws.Range("I28").Value = 170333.33
$ws.Range("J28").Value = 173333.33
$ws.Range("K28").Value = 170333.33
$ws.Range("L28").Value = 173333.33
$ws.Range("M28").Value = -170101.33
$ws.Range("N28").Value = -173797.33
$ws.Range("H29").Value = 16016
$ws.Range("I29").Value = 16016
$ws.Range("K29").Value = 16016
$ws.Range("M29").Value = -15721
$ws.Range("H31").Value = 677.3333
$ws.Range("I31").Value = 515
$ws.Range("J31").Value = 758.5
$ws.Range("K31").Value = 515
$ws.Range("L31").Value = 758.5
$ws.Range("M31").Value = -267
$ws.Range("N31").Value = -1254.5
$ws.Range("H37").Value = 171333.33
$ws.Range("I37").Value = 170333.33
$ws.Range("J37").Value = 173333.33
$ws.Range("K37").Value = 170333.33
$ws.Range("L37").Value = 173333.33
$ws.Range("M37").Value = -170226.33
$ws.Range("N37").Value = -173547.33
$ws.Range("H46").Value = 35844
$ws.Range("J46").Value = 3343.1428
$ws.Range("L46").Value = 3343.1428
$ws.Range("N46").Value = -3719.1428
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 2628.1667
$ws.Range("I61").Value = 2599.6667
$ws.Range("K61").Value = 2599.6667
$ws.Range("M61").Value = -2397.6667
$ws.Range("H63").Value = 74994.664
$ws.Range("J63").Value = 74994.664
$ws.Range("L63").Value = 74994.664
$ws.Range("N63").Value = -76492.664
$ws.Range("H66").Value = 74994.664
$ws.Range("J66").Value = 74994.664
$ws.Range("L66").Value = 224983.992
$ws.Range("N66").Value = -232471.992
$ws.Range("H74").Value = 49216.5
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 49216.5
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H100").Value = 2874
$ws.Range("I100").Value = 2950.4
$ws.Range("J100").Value = 2746.6667
$ws.Range("K100").Value = 2950.4
$ws.Range("L100").Value = 2746.6667
$ws.Range("M100").Value = -2409.4
$ws.Range("N100").Value = -3828.6667
$ws.Range("H113").Value = 2628.1667
$ws.Range("I113").Value = 2599.6667
$ws.Range("K113").Value = 2599.6667
$ws.Range("M113").Value = -429.6667000000002
$ws.Range("H136").Value = 3994.3572
$ws.Range("I136").Value = 3366.9092
$ws.Range("K136").Value = 10100.7276
$ws.Range("M136").Value = -7550.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16875.223
$ws.Range("J62").Value = 14125
$ws.Range("L62").Value = 14125
$ws.Range("N62").Value = -15373
$ws.Range("H65").Value = 16875.223
$ws.Range("J65").Value = 14125
$ws.Range("L65").Value = 70625
$ws.Range("N65").Value = -76865
$ws.Range("H70").Value = 21069.666
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 21069.666
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 1113072.1
$ws.Range("I81").Value = 2205.875
$ws.Range("K81").Value = 4411.75
$ws.Range("M81").Value = -3350.75
$ws.Range("H84").Value = 1113072.1
$ws.Range("I84").Value = 2205.875
$ws.Range("K84").Value = 22058.75
$ws.Range("M84").Value = -16754.75
$ws.Range("H122").Value = 1074.75
$ws.Range("I122").Value = 1074.75
$ws.Range("K122").Value = 3224.25
$ws.Range("M122").Value = -774.25
$ws.Range("H132").Value = 2377.72
$ws.Range("I132").Value = 2613.158
$ws.Range("J132").Value = 1632.1666
$ws.Range("K132").Value = 7839.474
$ws.Range("L132").Value = 4896.4998
$ws.Range("M132").Value = -5309.474
$ws.Range("N132").Value = -9956.4998
$ws.Range("H136").Value = 2854.2083
$ws.Range("J136").Value = 1413.5
$ws.Range("L136").Value = 4240.5
$ws.Range("N136").Value = -9340.5
